$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1296.8462
$ws.Range("I86").Value = 1325.8334
$ws.Range("J86").Value = 1272
$ws.Range("K86").Value = 1325.8334
$ws.Range("L86").Value = 1272
$ws.Range("M86").Value = -202.8334
$ws.Range("N86").Value = -3518
$ws.Range("H89").Value = 1296.8462
$ws.Range("I89").Value = 1325.8334
$ws.Range("J89").Value = 1272
$ws.Range("K89").Value = 6629.166999999999
$ws.Range("L89").Value = 6360
$ws.Range("M89").Value = -1013.166999999999
$ws.Range("N89").Value = -17592
$ws.Range("H137").Value = 759058.0600000001
$ws.Range("I137").Value = 2074235.8
$ws.Range("J137").Value = 2830.825
$ws.Range("K137").Value = 6222707.4
$ws.Range("L137").Value = 8492.474999999999
$ws.Range("M137").Value = -6220157.4
$ws.Range("N137").Value = -13592.475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 794.3333
$ws.Range("I110").Value = 667.4545000000001
$ws.Range("J110").Value = 993.7143
$ws.Range("K110").Value = 667.4545000000001
$ws.Range("L110").Value = 993.7143
$ws.Range("M110").Value = 1377.5455
$ws.Range("N110").Value = -5083.7143
$ws.Range("H137").Value = 40001
$ws.Range("J137").Value = 40001
$ws.Range("L137").Value = 40001
$ws.Range("N137").Value = -50201

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 278330.5
$ws.Range("I31").Value = 712184.8
$ws.Range("J31").Value = 3556.1333
$ws.Range("K31").Value = 712184.8
$ws.Range("L31").Value = 3556.1333
$ws.Range("M31").Value = -711889.8
$ws.Range("N31").Value = -4146.1333
$ws.Range("H34").Value = 278330.5
$ws.Range("I34").Value = 712184.8
$ws.Range("J34").Value = 3556.1333
$ws.Range("K34").Value = 712184.8
$ws.Range("L34").Value = 3556.1333
$ws.Range("M34").Value = -711982.8
$ws.Range("N34").Value = -3960.1333
$ws.Range("H60").Value = 28797.666
$ws.Range("J60").Value = 28797.666
$ws.Range("L60").Value = 28797.666
$ws.Range("N60").Value = -29819.666
$ws.Range("H74").Value = 32536.625
$ws.Range("J74").Value = 35858.285
$ws.Range("L74").Value = 35858.285
$ws.Range("N74").Value = -37606.285
$ws.Range("H77").Value = 32536.625
$ws.Range("J77").Value = 35858.285
$ws.Range("L77").Value = 107574.855
$ws.Range("N77").Value = -116310.855
$ws.Range("H105").Value = 1227.1
$ws.Range("I105").Value = 1139.1428
$ws.Range("J105").Value = 1432.3334
$ws.Range("K105").Value = 1139.1428
$ws.Range("L105").Value = 1432.3334
$ws.Range("M105").Value = 607.8571999999999
$ws.Range("N105").Value = -4926.3334
$ws.Range("H132").Value = 4671.778
$ws.Range("I132").Value = 3567.353
$ws.Range("K132").Value = 10702.059
$ws.Range("M132").Value = -8172.059000000001
$ws.Range("H137").Value = 44608.57
$ws.Range("J137").Value = 44608.57
$ws.Range("L137").Value = 44608.57
$ws.Range("N137").Value = -54808.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 22332
$ws.Range("J39").Value = 22332
$ws.Range("L39").Value = 66996
$ws.Range("N39").Value = -67584
$ws.Range("H40").Value = 278.1875
$ws.Range("I40").Value = 96.73333
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 386.93332
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = -317.93332
$ws.Range("N40").Value = -12138
$ws.Range("H58").Value = 4218.5
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 4714
$ws.Range("K58").Value = 2250
$ws.Range("L58").Value = 14142
$ws.Range("N58").Value = -14398
$ws.Range("M58").Value = -2122
$ws.Range("H64").Value = 1757.1428
$ws.Range("I64").Value = 480
$ws.Range("J64").Value = 4950
$ws.Range("K64").Value = 1440
$ws.Range("L64").Value = 14850
$ws.Range("M64").Value = -1170
$ws.Range("N64").Value = -15390
$ws.Range("H67").Value = 1757.1428
$ws.Range("I67").Value = 480
$ws.Range("J67").Value = 4950
$ws.Range("K67").Value = 1440
$ws.Range("L67").Value = 14850
$ws.Range("M67").Value = -504
$ws.Range("N67").Value = -16722
$ws.Range("H69").Value = 3407.7222
$ws.Range("I69").Value = 797.1429000000001
$ws.Range("J69").Value = 5069
$ws.Range("K69").Value = 2391.4287
$ws.Range("L69").Value = 15207
$ws.Range("M69").Value = -1580.4287
$ws.Range("N69").Value = -16829
$ws.Range("H72").Value = 3407.7222
$ws.Range("I72").Value = 797.1429000000001
$ws.Range("J72").Value = 5069
$ws.Range("K72").Value = 7174.2861
$ws.Range("L72").Value = 45621
$ws.Range("M72").Value = -3118.2861
$ws.Range("N72").Value = -53733
$ws.Range("H81").Value = 1418.8334
$ws.Range("J81").Value = 2500
$ws.Range("L81").Value = 7500
$ws.Range("N81").Value = -9746
$ws.Range("H84").Value = 1418.8334
$ws.Range("J84").Value = 2500
$ws.Range("L84").Value = 22500
$ws.Range("N84").Value = -33732
$ws.Range("H86").Value = 870.6896400000001
$ws.Range("J86").Value = 1137.5
$ws.Range("L86").Value = 3412.5
$ws.Range("N86").Value = -5784.5
$ws.Range("H89").Value = 870.6896400000001
$ws.Range("J89").Value = 1137.5
$ws.Range("L89").Value = 10237.5
$ws.Range("N89").Value = -22093.5
$ws.Range("H94").Value = 3826.6667
$ws.Range("J94").Value = 3826.6667
$ws.Range("L94").Value = 11480.0001
$ws.Range("N94").Value = -12832.0001
$ws.Range("H107").Value = 14505.689
$ws.Range("I107").Value = 410.56757
$ws.Range("J107").Value = 28600.81
$ws.Range("K107").Value = 1231.70271
$ws.Range("L107").Value = 85802.43000000001
$ws.Range("M107").Value = 688.29729
$ws.Range("N107").Value = -89642.43000000001
$ws.Range("H109").Value = 3759.5
$ws.Range("I109").Value = 613.5
$ws.Range("J109").Value = 5332.5
$ws.Range("K109").Value = 1840.5
$ws.Range("L109").Value = 15997.5
$ws.Range("M109").Value = -800.5
$ws.Range("N109").Value = -18077.5
$ws.Range("H113").Value = 4808511
$ws.Range("I113").Value = 670.2857
$ws.Range("J113").Value = 10417658
$ws.Range("K113").Value = 2010.8571
$ws.Range("L113").Value = 31252974
$ws.Range("M113").Value = 159.1428999999998
$ws.Range("N113").Value = -31257314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1352.9
$ws.Range("I97").Value = 1507.4
$ws.Range("J97").Value = 1198.4
$ws.Range("K97").Value = 1507.4
$ws.Range("L97").Value = 1198.4
$ws.Range("M97").Value = -1011.4
$ws.Range("N97").Value = -2190.4
$ws.Range("H122").Value = 4978.4546
$ws.Range("I122").Value = 3200.6
$ws.Range("K122").Value = 9601.799999999999
$ws.Range("M122").Value = -7151.799999999999
$ws.Range("H132").Value = 3870.5652
$ws.Range("I132").Value = 1862.2
$ws.Range("J132").Value = 4428.4443
$ws.Range("K132").Value = 5586.6
$ws.Range("L132").Value = 13285.3329
$ws.Range("M132").Value = -3056.6
$ws.Range("N132").Value = -18345.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1357.875
$ws.Range("I61").Value = 1422.5714
$ws.Range("J61").Value = 905
$ws.Range("K61").Value = 1422.5714
$ws.Range("L61").Value = 905
$ws.Range("M61").Value = -1220.5714
$ws.Range("N61").Value = -1309
$ws.Range("H113").Value = 1357.875
$ws.Range("I113").Value = 1422.5714
$ws.Range("J113").Value = 905
$ws.Range("K113").Value = 1422.5714
$ws.Range("L113").Value = 905
$ws.Range("M113").Value = 747.4286
$ws.Range("N113").Value = -5245

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 30000
$ws.Range("J26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("N26").Value = -30586
$ws.Range("H107").Value = 836.6875
$ws.Range("I107").Value = 601.1667
$ws.Range("K107").Value = 1803.5001
$ws.Range("M107").Value = 116.4999
$ws.Range("H122").Value = 3546.875
$ws.Range("I122").Value = 2184
$ws.Range("J122").Value = 4479.3687
$ws.Range("K122").Value = 6552
$ws.Range("L122").Value = 13438.1061
$ws.Range("M122").Value = -4102
$ws.Range("N122").Value = -18338.1061
